$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6920.2104
$ws.Range("J17").Value = 1701.8379
$ws.Range("L17").Value = 5105.5137
$ws.Range("N17").Value = -5441.5137

$ws.Range("H44").Value = 8000
$ws.Range("J44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("N44").Value = -8924

$ws.Range("H92").Value = 775.55554
$ws.Range("I92").Value = 736
$ws.Range("J92").Value = 825
$ws.Range("K92").Value = 736
$ws.Range("L92").Value = 825
$ws.Range("M92").Value = 512
$ws.Range("N92").Value = -3321

$ws.Range("H100").Value = 3385.5557
$ws.Range("I100").Value = 2932.8333
$ws.Range("J100").Value = 4291
$ws.Range("K100").Value = 2932.8333
$ws.Range("L100").Value = 4291
$ws.Range("M100").Value = -2391.8333
$ws.Range("N100").Value = -5373

$ws.Range("H135").Value = 1019.4375
$ws.Range("I135").Value = 820.73334
$ws.Range("K135").Value = 7386.60006
$ws.Range("M135").Value = -4851.60006

$ws.Range("H138").Value = 4390.533
$ws.Range("I138").Value = 2981.2727
$ws.Range("J138").Value = 4632.75
$ws.Range("K138").Value = 8943.8181
$ws.Range("L138").Value = 13898.25
$ws.Range("M138").Value = -3803.8181
$ws.Range("N138").Value = -24178.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 21156
$ws.Range("I41").Value = 1593.6
$ws.Range("K41").Value = 1593.6
$ws.Range("M41").Value = -1179.6

$ws.Range("H122").Value = 2611.8948
$ws.Range("I122").Value = 1793.2307
$ws.Range("J122").Value = 4385.6665
$ws.Range("K122").Value = 5379.6921
$ws.Range("L122").Value = 13156.9995
$ws.Range("M122").Value = -2929.6921
$ws.Range("N122").Value = -18056.9995

$ws.Range("H132").Value = 12660103
$ws.Range("I132").Value = 17242348
$ws.Range("J132").Value = 4380.5713
$ws.Range("K132").Value = 51727044
$ws.Range("L132").Value = 13141.7139
$ws.Range("M132").Value = -51724514
$ws.Range("N132").Value = -18201.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 763.375
$ws.Range("I64").Value = 600
$ws.Range("J64").Value = 861.4
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 861.4
$ws.Range("M64").Value = -375
$ws.Range("N64").Value = -1311.4

$ws.Range("H67").Value = 763.375
$ws.Range("I67").Value = 600
$ws.Range("J67").Value = 861.4
$ws.Range("K67").Value = 600
$ws.Range("L67").Value = 861.4
$ws.Range("M67").Value = 180
$ws.Range("N67").Value = -2421.4

$ws.Range("H99").Value = 3440
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 4850
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 4850
$ws.Range("M99").Value = -1002
$ws.Range("N99").Value = -7846

$ws.Range("H134").Value = 3535.4614
$ws.Range("I134").Value = 2433.875
$ws.Range("J134").Value = 5298
$ws.Range("K134").Value = 7301.625
$ws.Range("L134").Value = 15894
$ws.Range("M134").Value = -4766.625
$ws.Range("N134").Value = -20964

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3675
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 4666.6665
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 4666.6665
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -5366.6665

$ws.Range("H74").Value = 17414.273
$ws.Range("J74").Value = 17414.273
$ws.Range("L74").Value = 17414.273
$ws.Range("N74").Value = -19162.273

$ws.Range("H77").Value = 17414.273
$ws.Range("J77").Value = 17414.273
$ws.Range("L77").Value = 52242.819
$ws.Range("N77").Value = -60978.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 7155
$ws.Range("I22").Value = 990.5
$ws.Range("J22").Value = 8524.888999999999
$ws.Range("K22").Value = 2971.5
$ws.Range("L22").Value = 25574.667
$ws.Range("M22").Value = -2802.5
$ws.Range("N22").Value = -25912.667

$ws.Range("H23").Value = 84
$ws.Range("I23").Value = 70
$ws.Range("J23").Value = 91
$ws.Range("K23").Value = 210
$ws.Range("L23").Value = 273
$ws.Range("M23").Value = 25
$ws.Range("N23").Value = -743

$ws.Range("H27").Value = 7155
$ws.Range("I27").Value = 990.5
$ws.Range("J27").Value = 8524.888999999999
$ws.Range("K27").Value = 2971.5
$ws.Range("L27").Value = 25574.667
$ws.Range("M27").Value = -2869.5
$ws.Range("N27").Value = -25778.667

$ws.Range("H40").Value = 155.2
$ws.Range("I40").Value = 107.42857
$ws.Range("J40").Value = 266.66666
$ws.Range("K40").Value = 429.71428
$ws.Range("L40").Value = 1066.66664
$ws.Range("M40").Value = -360.71428
$ws.Range("N40").Value = -1204.66664

$ws.Range("H61").Value = 10244.917
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 10244.917
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30734.751
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -31164.751

$ws.Range("H131").Value = 1197.1389
$ws.Range("I131").Value = 1229.6923
$ws.Range("J131").Value = 1178.7391
$ws.Range("K131").Value = 3689.0769
$ws.Range("L131").Value = 3536.2173
$ws.Range("M131").Value = 1350.9231
$ws.Range("N131").Value = -13616.2173

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 65099.812
$ws.Range("I102").Value = 1899.8
$ws.Range("J102").Value = 170433.17
$ws.Range("K102").Value = 1899.8
$ws.Range("L102").Value = 170433.17
$ws.Range("M102").Value = -277.8
$ws.Range("N102").Value = -173677.17

$ws.Range("H132").Value = 2908.6938
$ws.Range("I132").Value = 2436.147
$ws.Range("J132").Value = 3979.8
$ws.Range("K132").Value = 7308.441
$ws.Range("L132").Value = 11939.4
$ws.Range("M132").Value = -4778.441
$ws.Range("N132").Value = -16999.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 464642.94
$ws.Range("I2").Value = 461923.06
$ws.Range("J2").Value = 500001
$ws.Range("K2").Value = 461923.06
$ws.Range("L2").Value = 500001
$ws.Range("M2").Value = -461811.06
$ws.Range("N2").Value = -500225

$ws.Range("H17").Value = 28888
$ws.Range("I17").Value = 28888
$ws.Range("K17").Value = 28888
$ws.Range("M17").Value = -28718

$ws.Range("H22").Value = 90911310
$ws.Range("I22").Value = 166667740
$ws.Range("J22").Value = 3580.4
$ws.Range("K22").Value = 166667740
$ws.Range("L22").Value = 3580.4
$ws.Range("M22").Value = -166667445
$ws.Range("N22").Value = -4170.4

$ws.Range("H27").Value = 90911310
$ws.Range("I27").Value = 166667740
$ws.Range("J27").Value = 3580.4
$ws.Range("K27").Value = 166667740
$ws.Range("L27").Value = 3580.4
$ws.Range("M27").Value = -166667633
$ws.Range("N27").Value = -3794.4

$ws.Range("H132").Value = 3352.2188
$ws.Range("I132").Value = 1876.7646
$ws.Range("J132").Value = 5024.4
$ws.Range("K132").Value = 5630.293799999999
$ws.Range("L132").Value = 15073.2
$ws.Range("M132").Value = -3100.293799999999
$ws.Range("N132").Value = -20133.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1697.7778
$ws.Range("I107").Value = 635.7
$ws.Range("J107").Value = 3025.375
$ws.Range("K107").Value = 1907.1
$ws.Range("L107").Value = 9076.125
$ws.Range("M107").Value = 12.89999999999986
$ws.Range("N107").Value = -12916.125

$ws.Range("H133").Value = 42245
$ws.Range("J133").Value = 42245
$ws.Range("L133").Value = 42245
$ws.Range("N133").Value = -52365
